# Update "想去人数" (want-to-go count) values in column F on the
# "展览" and "全部类型" worksheets, matching the refreshed data export.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetAllTypes   = $wb.Worksheets.Item("全部类型")

# Mapping of row -> new value for the "展览" sheet
$exhibitionUpdates = @{
    2  = 2234
    3  = 105
    4  = 13551
    8  = 496
    9  = 1202
    10 = 1017
    11 = 13855
    12 = 14636
    16 = 46
    20 = 13
    21 = 54
    23 = 1130
    26 = 5627
    28 = 1046
    29 = 5377
    30 = 40
    31 = 38
    32 = 198
}

foreach ($row in $exhibitionUpdates.Keys) {
    $sheetExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

# Mapping of row -> new value for the "全部类型" sheet
$allTypesUpdates = @{
    2  = 2234
    3  = 105
    4  = 13551
    9  = 496
    10 = 1202
    11 = 1017
    12 = 13855
    13 = 14636
    17 = 46
    21 = 13
    22 = 54
    23 = 10
    24 = 1130
    27 = 5627
    29 = 1046
    30 = 5377
    31 = 40
    32 = 38
    33 = 198
}

foreach ($row in $allTypesUpdates.Keys) {
    $sheetAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}

$wb.Save()
